# Apply updated TPM-derived values (ligand/receptor expression + derived specificity/edge weights)
# for the Il15 -> Il2ra LR-pair sheet. Values correspond to the refreshed TPM input data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 7.286709999999999
    "H2" = 21.86013
    "I2" = 0.1632739668438106
    "J2" = 0.1632739668438107
    "M2" = 0.01512566666666667
    "N2" = 0.045377
    "O2" = 0.01277774711566512
    "P2" = 0.01277774711566512
    "Q2" = 0.1102163465566667
    "R2" = 0.9919471190099999
    "S2" = 0.002086273458901704
    "T2" = 0.002086273458901705
    "G3" = 7.286709999999999
    "H3" = 21.86013
    "I3" = 0.1632739668438106
    "J3" = 0.1632739668438107
    "M3" = 0.8534100000000001
    "N3" = 2.56023
    "O3" = 0.7209372919747741
    "P3" = 0.7209372919747741
    "Q3" = 6.2185511811
    "R3" = 55.9669606299
    "S3" = 0.1177102915063559
    "T3" = 0.1177102915063559
    "G4" = 7.286709999999999
    "H4" = 21.86013
    "I4" = 0.1632739668438106
    "J4" = 0.1632739668438107
    "M4" = 0.1323446666666667
    "N4" = 0.397034
    "O4" = 0.1118011337973199
    "P4" = 0.1118011337973199
    "Q4" = 0.9643572060466665
    "R4" = 8.67921485442
    "S4" = 0.01825421461272405
    "T4" = 0.01825421461272405
    "G5" = 7.286709999999999
    "H5" = 21.86013
    "I5" = 0.1632739668438106
    "J5" = 0.1632739668438107
    "M5" = 0.1828703333333333
    "N5" = 0.548611
    "O5" = 0.154483827112241
    "P5" = 0.154483827112241
    "Q5" = 1.332523086603333
    "R5" = 11.99270777943
    "S5" = 0.02522318726582901
    "T5" = 0.02522318726582901
    "I6" = 0.05961564640488622
    "J6" = 0.05961564640488623
    "M6" = 0.01512566666666667
    "N6" = 0.045377
    "O6" = 0.01277774711566512
    "P6" = 0.01277774711566512
    "Q6" = 0.04024290504711112
    "R6" = 0.362186145424
    "S6" = 0.0007617536538985469
    "T6" = 0.000761753653898547
    "I7" = 0.05961564640488622
    "J7" = 0.05961564640488623
    "M7" = 0.8534100000000001
    "N7" = 2.56023
    "O7" = 0.7209372919747741
    "P7" = 0.7209372919747741
    "Q7" = 2.27055761264
    "R7" = 20.43501851376
    "S7" = 0.04297914267846435
    "T7" = 0.04297914267846435
    "I8" = 0.05961564640488622
    "J8" = 0.05961564640488623
    "M8" = 0.1323446666666667
    "N8" = 0.397034
    "O8" = 0.1118011337973199
    "P8" = 0.1118011337973199
    "Q8" = 0.3521123380231111
    "R8" = 3.169011042208
    "S8" = 0.006665096860126399
    "T8" = 0.0066650968601264
    "I9" = 0.05961564640488622
    "J9" = 0.05961564640488623
    "M9" = 0.1828703333333333
    "N9" = 0.548611
    "O9" = 0.154483827112241
    "P9" = 0.154483827112241
    "Q9" = 0.4865394446702222
    "R9" = 4.378855002032
    "S9" = 0.009209653212396933
    "T9" = 0.009209653212396935
    "G10" = 19.00851733333333
    "H10" = 57.025552
    "I10" = 0.4259255588369328
    "J10" = 0.4259255588369329
    "M10" = 0.01512566666666667
    "N10" = 0.045377
    "O10" = 0.01277774711566512
    "P10" = 0.01277774711566512
    "Q10" = 0.2875164970115556
    "R10" = 2.587648473104
    "S10" = 0.005442369080916674
    "T10" = 0.005442369080916675
    "G11" = 19.00851733333333
    "H11" = 57.025552
    "I11" = 0.4259255588369328
    "J11" = 0.4259255588369329
    "M11" = 0.8534100000000001
    "N11" = 2.56023
    "O11" = 0.7209372919747741
    "P11" = 0.7209372919747741
    "Q11" = 16.22205877744
    "R11" = 145.99852899696
    "S11" = 0.3070656189707406
    "T11" = 0.3070656189707407
    "G12" = 19.00851733333333
    "H12" = 57.025552
    "I12" = 0.4259255588369328
    "J12" = 0.4259255588369329
    "M12" = 0.1323446666666667
    "N12" = 0.397034
    "O12" = 0.1118011337973199
    "P12" = 0.1118011337973199
    "Q12" = 2.515675890307556
    "R12" = 22.641083012768
    "S12" = 0.04761896039122619
    "T12" = 0.0476189603912262
    "G13" = 19.00851733333333
    "H13" = 57.025552
    "I13" = 0.4259255588369328
    "J13" = 0.4259255588369329
    "M13" = 0.1828703333333333
    "N13" = 0.548611
    "O13" = 0.154483827112241
    "P13" = 0.154483827112241
    "Q13" = 3.476093900919111
    "R13" = 31.284845108272
    "S13" = 0.06579861039404936
    "T13" = 0.06579861039404937
    "G14" = 0.1146546666666667
    "H14" = 0.343964
    "I14" = 0.002569077436020028
    "J14" = 0.002569077436020028
    "M14" = 0.01512566666666667
    "N14" = 0.045377
    "O14" = 0.01277774711566512
    "P14" = 0.01277774711566512
    "Q14" = 0.001734228269777778
    "R14" = 0.015608054428
    "S14" = 0.00003282702179802526
    "T14" = 0.00003282702179802526
    "G15" = 0.1146546666666667
    "H15" = 0.343964
    "I15" = 0.002569077436020028
    "J15" = 0.002569077436020028
    "M15" = 0.8534100000000001
    "N15" = 2.56023
    "O15" = 0.7209372919747741
    "P15" = 0.7209372919747741
    "Q15" = 0.09784743908000001
    "R15" = 0.8806269517200001
    "S15" = 0.001852143729597775
    "T15" = 0.001852143729597775
    "G16" = 0.1146546666666667
    "H16" = 0.343964
    "I16" = 0.002569077436020028
    "J16" = 0.002569077436020028
    "M16" = 0.1323446666666667
    "N16" = 0.397034
    "O16" = 0.1118011337973199
    "P16" = 0.1118011337973199
    "Q16" = 0.01517393364177778
    "R16" = 0.136565402776
    "S16" = 0.0002872257701601508
    "T16" = 0.0002872257701601508
    "G17" = 0.1146546666666667
    "H17" = 0.343964
    "I17" = 0.002569077436020028
    "J17" = 0.002569077436020028
    "M17" = 0.1828703333333333
    "N17" = 0.548611
    "O17" = 0.154483827112241
    "P17" = 0.154483827112241
    "Q17" = 0.02096693711155555
    "R17" = 0.188702434004
    "S17" = 0.0003968809144640773
    "T17" = 0.0003968809144640773
    "G18" = 15.55827866666667
    "H18" = 46.674836
    "I18" = 0.3486157504783503
    "J18" = 0.3486157504783503
    "M18" = 0.01512566666666667
    "N18" = 0.045377
    "O18" = 0.01277774711566512
    "P18" = 0.01277774711566512
    "Q18" = 0.2353293370191111
    "R18" = 2.117964033172
    "S18" = 0.004454523900150172
    "T18" = 0.004454523900150173
    "G19" = 15.55827866666667
    "H19" = 46.674836
    "I19" = 0.3486157504783503
    "J19" = 0.3486157504783503
    "M19" = 0.8534100000000001
    "N19" = 2.56023
    "O19" = 0.7209372919747741
    "P19" = 0.7209372919747741
    "Q19" = 13.27759059692
    "R19" = 119.49831537228
    "S19" = 0.2513300950896154
    "T19" = 0.2513300950896154
    "G20" = 15.55827866666667
    "H20" = 46.674836
    "I20" = 0.3486157504783503
    "J20" = 0.3486157504783503
    "M20" = 0.1323446666666667
    "N20" = 0.397034
    "O20" = 0.1118011337973199
    "P20" = 0.1118011337973199
    "Q20" = 2.059055204047111
    "R20" = 18.531496836424
    "S20" = 0.03897563616308314
    "T20" = 0.03897563616308314
    "G21" = 15.55827866666667
    "H21" = 46.674836
    "I21" = 0.3486157504783503
    "J21" = 0.3486157504783503
    "M21" = 0.1828703333333333
    "N21" = 0.548611
    "O21" = 0.154483827112241
    "P21" = 0.154483827112241
    "Q21" = 2.845147605866222
    "R21" = 25.606328452796
    "S21" = 0.0538554953255016
    "T21" = 0.05385549532550161
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

